# The workbook's header row (row 1) is being relabeled to reflect the new
# "reference levels" naming scheme used by the newest version of the demo.
#
#   B1: min       -> pess
#   C1: SD_nedre  -> X0
#   D1: SD_D      -> X20
#   E1: D_M       -> X40
#   F1: M_G       -> X60
#   G1: G_SG      -> X80
#   H1: SG_ovre   -> X100
#   I1: max       -> opt
#
# Column A (the "typ" labels CS1..CB5) and all numeric class-limit data in
# columns B:I for rows 2-30 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "pess"
$ws.Range("C1").Value = "X0"
$ws.Range("D1").Value = "X20"
$ws.Range("E1").Value = "X40"
$ws.Range("F1").Value = "X60"
$ws.Range("G1").Value = "X80"
$ws.Range("H1").Value = "X100"
$ws.Range("I1").Value = "opt"
